$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("G3").Value = 3.4
$ws.Range("I3").Value = 2.45
$ws.Range("O3").Value = 1.73
$ws.Range("P3").Value = 2
$ws.Range("S3").Value = 1.75
$ws.Range("AY3").Value = 41
$ws.Range("BA3").Value = 126
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 2.7
$ws.Range("R4").Value = 1.44
$ws.Range("R5").Value = 1.91
$ws.Range("M9").Value = 1.05
$ws.Range("O9").Value = 1.33
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 1.66
$ws.Range("M10").Value = 1.05
$ws.Range("O10").Value = 1.33
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 1.62
$ws.Range("X10").Value = 13
$ws.Range("G11").Value = 2.5
$ws.Range("I11").Value = 2.8
$ws.Range("J11").Value = 3.2
$ws.Range("K11").Value = 2.05
$ws.Range("L11").Value = 3.5
$ws.Range("M11").Value = 1.05
$ws.Range("O11").Value = 1.3
$ws.Range("Q11").Value = 2.03
$ws.Range("R11").Value = 1.68
$ws.Range("Z11").Value = 23
$ws.Range("AP11").Value = 23
$ws.Range("AX11").Value = 17
$ws.Range("J12").Value = 2.4
$ws.Range("L12").Value = 5
$ws.Range("M12").Value = 1.04
$ws.Range("O12").Value = 1.27
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 1.8
$ws.Range("U12").Value = 1.91
$ws.Range("V12").Value = 1.91
$ws.Range("X12").Value = 8
$ws.Range("AC12").Value = 9.5
$ws.Range("AX12").Value = 26
$ws.Range("I13").Value = 2.75
$ws.Range("U13").Value = 2.1
$ws.Range("V13").Value = 1.67
$ws.Range("AZ13").Value = 51
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.8
$ws.Range("G16").Value = 1.5
$ws.Range("AG16").Value = 21
$ws.Range("AN16").Value = 3.6
$ws.Range("G17").Value = 1.41
$ws.Range("M17").Value = 1.03
$ws.Range("O17").Value = 1.22
$ws.Range("G18").Value = 2.15
$ws.Range("I18").Value = 3.3
$ws.Range("J18").Value = 2.88
$ws.Range("K18").Value = 2.1
$ws.Range("M18").Value = 1.05
$ws.Range("O18").Value = 1.33
$ws.Range("Y18").Value = 9
$ws.Range("Z18").Value = 19
$ws.Range("AB18").Value = 29
$ws.Range("AC18").Value = 9
$ws.Range("AE18").Value = 15
$ws.Range("AG18").Value = 9
$ws.Range("AH18").Value = 17
$ws.Range("AM18").Value = 301
$ws.Range("AN18").Value = 4
$ws.Range("AO18").Value = 12
$ws.Range("AS18").Value = 151
$ws.Range("AT18").Value = 2.62
$ws.Range("AU18").Value = 8
$ws.Range("BB18").Value = 201
$ws.Range("M19").Value = 1.04
$ws.Range("O19").Value = 1.27
$ws.Range("M20").Value = 1.05
$ws.Range("N20").Value = 8
$ws.Range("O20").Value = 1.37
$ws.Range("BD20").Value = 151
$ws.Range("G21").Value = 1.71
$ws.Range("H21").Value = 3.7
$ws.Range("I21").Value = 4.33
$ws.Range("J21").Value = 2.3
$ws.Range("L21").Value = 4.5
$ws.Range("M21").Value = 1.02
$ws.Range("O21").Value = 1.15
$ws.Range("U21").Value = 1.62
$ws.Range("V21").Value = 2.2
$ws.Range("W21").Value = 9
$ws.Range("X21").Value = 9.5
$ws.Range("Z21").Value = 15
$ws.Range("AD21").Value = 7.5
$ws.Range("AE21").Value = 13
$ws.Range("AI21").Value = 15
$ws.Range("AK21").Value = 34
$ws.Range("AM21").Value = 151
$ws.Range("AO21").Value = 9
$ws.Range("AY21").Value = 26
$ws.Range("Q22").Value = 1.77
$ws.Range("R22").Value = 1.97
$ws.Range("G23").Value = 3.2
$ws.Range("H23").Value = 3.3
$ws.Range("I23").Value = 2.2
$ws.Range("L23").Value = 2.87
$ws.Range("N23").Value = 10
$ws.Range("O23").Value = 1.29
$ws.Range("P23").Value = 3.5
$ws.Range("Q23").Value = 1.94
$ws.Range("R23").Value = 1.79
$ws.Range("AC23").Value = 10
$ws.Range("AD23").Value = 6.5
$ws.Range("AF23").Value = 41
$ws.Range("AK23").Value = 17
$ws.Range("AL23").Value = 26
$ws.Range("AY23").Value = 21
$ws.Range("BA23").Value = 51
$ws.Range("N24").Value = 13
$ws.Range("O24").Value = 1.22
$ws.Range("P24").Value = 4
$ws.Range("Q24").Value = 1.8
$ws.Range("BD24").Value = 151
$ws.Range("Q25").Value = 1.87
$ws.Range("R25").Value = 1.87
$ws.Range("G29").Value = 1.7
$ws.Range("H29").Value = 3.9
$ws.Range("I29").Value = 4.33
$ws.Range("J29").Value = 2.3
$ws.Range("AG29").Value = 13
$ws.Range("AO29").Value = 9
$ws.Range("AP29").Value = 19
$ws.Range("AQ29").Value = 29
$ws.Range("BB29").Value = 201
$ws.Range("G34").Value = 2.75
$ws.Range("H34").Value = 3.2
$ws.Range("I34").Value = 2.55
$ws.Range("J34").Value = 3.4
$ws.Range("L34").Value = 3.2
$ws.Range("M34").Value = 1.06
$ws.Range("N34").Value = 10
$ws.Range("O34").Value = 1.29
$ws.Range("P34").Value = 3.75
$ws.Range("Q34").Value = 1.93
$ws.Range("R34").Value = 1.93
$ws.Range("U34").Value = 1.7
$ws.Range("V34").Value = 2.05
$ws.Range("X34").Value = 15
$ws.Range("Z34").Value = 29
$ws.Range("AB34").Value = 29
$ws.Range("AE34").Value = 13
$ws.Range("AG34").Value = 9
$ws.Range("AH34").Value = 13
$ws.Range("AJ34").Value = 26
$ws.Range("AN34").Value = 4.75
$ws.Range("AO34").Value = 15
$ws.Range("AP34").Value = 23
$ws.Range("AQ34").Value = 51
$ws.Range("AW34").Value = 4.5
$ws.Range("AZ34").Value = 41
$ws.Range("M35").Value = 1.03
$ws.Range("N35").Value = 15
$ws.Range("Q35").Value = 1.63
$ws.Range("R36").Value = 1.62
$ws.Range("G37").Value = 2.6
$ws.Range("H37").Value = 3.25
$ws.Range("I37").Value = 2.63
$ws.Range("J37").Value = 3.25
$ws.Range("L37").Value = 3.25
$ws.Range("N37").Value = 10
$ws.Range("W37").Value = 8.5
$ws.Range("Y37").Value = 10
$ws.Range("AC37").Value = 10
$ws.Range("AE37").Value = 15
$ws.Range("AH37").Value = 13
$ws.Range("AN37").Value = 4.5
$ws.Range("AW37").Value = 4.75
